{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Target change (per the diff):\n//   1. The existing \"test-functional-spec\" paragraph gets wrapped with\n//      grammar proof-error markers (<w:proofErr w:type=\"gramStart\"/> ...\n//      <w:proofErr w:type=\"gramEnd\"/>).\n//   2. A brand-new paragraph is appended:\n//        \"editing 1st time \u2013 adding blah blah blah \u2026\"\n//      where \"st\" is superscript, and each \"blah\" word is wrapped in\n//      spell-check proof-error markers (spellStart/spellEnd).\n//   3. The \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd id=0) moves from\n//      the end of the first paragraph to the end of the new paragraph.\n//\n// <w:proofErr> markers aren't exposed as first-class Office.js objects, so\n// we build the exact target body markup and drop it in with\n// `body.insertOoxml(..., \"Replace\")`, which lets us control raw OOXML\n// (including proofErr and the bookmark) while leaving the rest of the\n// package (styles, sectPr, etc.) untouched.\n\nconst body = context.document.body;\n\nconst flatOpc = [\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>',\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">',\n  '  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">',\n  '    <pkg:xmlData>',\n  '      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">',\n  '        <w:body>',\n  '          <w:p>',\n  '            <w:proofErr w:type=\"gramStart\"/>',\n  '            <w:r><w:t>test-functional-spec</w:t></w:r>',\n  '            <w:proofErr w:type=\"gramEnd\"/>',\n  '          </w:p>',\n  '          <w:p>',\n  '            <w:r><w:t>editing 1</w:t></w:r>',\n  '            <w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>st</w:t></w:r>',\n  '            <w:r><w:t xml:space=\"preserve\"> time \\u2013 adding blah </w:t></w:r>',\n  '            <w:proofErr w:type=\"spellStart\"/>',\n  '            <w:r><w:t>blah</w:t></w:r>',\n  '            <w:proofErr w:type=\"spellEnd\"/>',\n  '            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>',\n  '            <w:proofErr w:type=\"spellStart\"/>',\n  '            <w:r><w:t>blah</w:t></w:r>',\n  '            <w:proofErr w:type=\"spellEnd\"/>',\n  '            <w:r><w:t xml:space=\"preserve\"> \\u2026</w:t></w:r>',\n  '            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>',\n  '            <w:bookmarkEnd w:id=\"0\"/>',\n  '          </w:p>',\n  '        </w:body>',\n  '      </w:document>',\n  '    </pkg:xmlData>',\n  '  </pkg:part>',\n  '</pkg:package>'\n].join('\\n');\n\nbody.insertOoxml(flatOpc, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# PowerShell-style Word COM interop edit script.\n#\n# Target change (per the diff):\n#   1. The existing \"test-functional-spec\" paragraph gets wrapped with\n#      grammar proof-error markers (<w:proofErr w:type=\"gramStart\"/> ...\n#      <w:proofErr w:type=\"gramEnd\"/>).\n#   2. A brand-new paragraph is appended:\n#        \"editing 1st time \u2013 adding blah blah blah \u2026\"\n#      where \"st\" is superscript, and each \"blah\" word is wrapped in\n#      spell-check proof-error markers (spellStart/spellEnd).\n#   3. The \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd id=0) moves from\n#      the end of the first paragraph to the end of the new paragraph.\n#\n# <w:proofErr> markers (and moving the bookmark) aren't reachable through\n# the high-level Paragraphs/Range.Text surface, so we drive Range.InsertXML\n# with the exact target body markup (flat-OPC wrapped), replacing the\n# document's whole Content range. That lets us control raw OOXML while\n# leaving the rest of the package (styles, sectPr, etc.) untouched.\n\n$d = $word.ActiveDocument\n$rng = $d.Content\n\n$xml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r><w:t>test-functional-spec</w:t></w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n          </w:p>\n          <w:p>\n            <w:r><w:t>editing 1</w:t></w:r>\n            <w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>st</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\"> time \u2013 adding blah </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>blah</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>blah</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> \u2026</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$rng.InsertXML($xml)\n"}
